$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Shorten the bullet text in "登录JMC developer web，先申请为开发者，然后创建
#    Application，再创建Agent..." -> remove "先" / "然后" / "再".
# ---------------------------------------------------------------------------
$oldBullet = "先申请为开发者，然后创建Application，再创建Agent"
$newBullet = "申请为开发者，创建Application，创建Agent"
$d.Content.Find.Execute($oldBullet, $true, $false, $false, $false, $false, $true, 1, $false, $newBullet, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge the two runs around the old "_GoBack" bookmark in the ChatApi
#    sentence back into a single run (this also removes that stray bookmark).
# ---------------------------------------------------------------------------
$chatApiSentence = "对于ChatApi消息，SDK会把收到的消息回调给第三方开发者的java demo，然后由第三方开发者调用sendChatCmd接口即可回复消息给发送者；"
$d.Content.Find.Execute($chatApiSentence, $true, $false, $false, $false, $false, $true, 1, $false, $chatApiSentence, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Delete the paragraphs between the "登录JMC developer web..." bullet and the
#    "实际开发流程" heading (the screenshots + captions are no longer needed).
# ---------------------------------------------------------------------------
$startRng = $d.Content
$startRng.Find.Execute("申请为开发者，创建Application，创建Agent，并将创建的Agent添加到创建的Application下，同时一并设置WebHookUrl；", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteFrom = $startRng.End + 1

$endRng = $d.Content
$endRng.Find.Execute("实际开发流程", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteTo = $endRng.Start

$killRange = $d.Range($deleteFrom, $deleteTo)
$killRange.Delete()

# ---------------------------------------------------------------------------
# 4) Bookmark "OLE_LINK1" around the URL prefix of the hyperlink display text.
# ---------------------------------------------------------------------------
$urlRng = $d.Content
$urlRng.Find.Execute("https://op-official.freepp.com/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("OLE_LINK1", $urlRng) | Out-Null

# ---------------------------------------------------------------------------
# 5) Re-create the "_GoBack" bookmark at the new edit point, mid-word inside
#    "WebHookUrl" (between "WebHo" and "okUrl"), matching where Word would
#    leave it after the bullet text was last edited.
# ---------------------------------------------------------------------------
$hookRng = $d.Content
$hookRng.Find.Execute("WebHookUrl", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $hookRng.Start + 5
$goBackRng = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $goBackRng) | Out-Null

# ---------------------------------------------------------------------------
# 6) Mark the "Hyperlink" character style as a Quick Style.
# ---------------------------------------------------------------------------
$hyperlinkStyle = $d.Styles.Item("Hyperlink")
$hyperlinkStyle.QuickStyle = $true
